$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '63.717.38'
$ws.Range("E2").Value = '  +1.95%  '

$ws.Range("D3").Value = '3.417.83'
$ws.Range("E3").Value = '  +2.26%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.48'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.04'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  +3.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = $plainStyle
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.418.23'
$ws.Range("E8").Value = '  +2.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.547'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  +1.94%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.39'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  -0.22%  '

$ws.Range("E11").Value = '  +4.27%  '

$ws.Range("E12").Value = '  -0.33%  '

$ws.Range("D13").Value = '4.010.80'
$ws.Range("E13").Value = '  +2.28%  '

$ws.Range("E14").Value = '  -3.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000193'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +8.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.22'
$ws.Range("D16").Style = $plainStyle
$ws.Range("E16").Value = '  +1.62%  '

$ws.Range("D17").Value = '63.735.10'
$ws.Range("E17").Value = '  +1.97%  '

$ws.Range("D18").Value = '3.414.90'
$ws.Range("E18").Value = '  +2.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.27'
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  -0.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("D20").Style = $plainStyle

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.28'
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = '  -1.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.07'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -4.13%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.65'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  +2.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.529'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("E26").Value = '  +26.99%  '

$ws.Range("E27").Value = '  +3.47%  '

$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.10'
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = '  +9.37%  '

$ws.Range("E31").Value = '  +1.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.36'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  +4.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.26'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  +1.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.38'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  -3.22%  '

$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.78'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  +0.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.80'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  +0.37%  '

$ws.Range("E38").Value = '  -1.54%  '

$ws.Range("D39").Value = '2.983.87'
$ws.Range("E39").Value = '  +7.09%  '

$ws.Range("E40").Value = '  +2.70%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  -2.00%  '

$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.02'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("E43").Value = '  -0.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.02'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  +3.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.761'
$ws.Range("D45").Style = $plainStyle
$ws.Range("E45").Value = '  +2.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.31'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  +1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.25'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  +5.63%  '

$ws.Range("E48").Value = '  +4.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.19'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  +23.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.832'
$ws.Range("D50").Style = $plainStyle
$ws.Range("E50").Value = '  +4.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.35'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  +0.71%  '

